$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert a new 2022-Q4 summary row above the existing
#    2022-Q3 row (which shifts down from row 2 to row 3).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Remember the existing 2022-Q3 totals before they get overwritten.
$q3Label  = $total.Range("B2").Value2
$q3Count  = $total.Range("C2").Value2
$q3Value  = $total.Range("D2").Value2

# Carry the index-column formatting (bold, centered, bordered) down to the
# row the old data is about to occupy.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

# Move 2022-Q3's totals down to row 3, renumbering its index.
$total.Range("A3").Value = 1
$total.Range("B3").Value = $q3Label
$total.Range("C3").Value = $q3Count
$total.Range("D3").Value = $q3Value

# Write the new 2022-Q4 totals into row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" sheet right before "2022-Q3" and fill it
#    with the quarter's fund holdings table.
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3Sheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "000646", "华润元大量化优选混合A", "0.19", "65.16", "0.61", "0.0012", 10),
    @(1, "007827", "华润元大量化优选混合C", "0.01", "65.16", "0.61", "0.0001", 10)
)

foreach ($row in $rows) {
    $r = $row[0] + 2
    $q4.Cells.Item($r, 1).Value = $row[0]

    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]

    # D/E/F/G hold numeric-looking figures that must stay text, matching the
    # other quarterly sheets (t="inlineStr").
    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 4).ClearFormats()

    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 5).ClearFormats()

    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 6).ClearFormats()

    $q4.Cells.Item($r, 7).NumberFormat = "@"
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 7).ClearFormats()

    $q4.Cells.Item($r, 8).Value = $row[7]
}

# Apply the bold/centered/bordered "header" style (matching the rest of the
# workbook) to the new sheet's header row and index column.
$total.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$total.Range("A2:A3").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)
